$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = 12345
$ws.Range("B2").Value = "Hesoyam"
$ws.Range("C2").Value = "SMK Negeri 1 Siatasbarita"
$ws.Range("D2").Value = 5000000
$ws.Range("E2").Value = "Teknologi Komputer"
$ws.Range("F2").Value = "UTBK"
$ws.Range("G2").Value = 2.4
$ws.Range("H2").Value = 3.1
$ws.Range("I2").Value = 3.2
$ws.Range("J2").Value = 3.5

# Row 3 updates
$ws.Range("A3").Value = 23456
$ws.Range("B3").Value = "Aezakmi"
$ws.Range("C3").Value = "SMA Negri 1 Sidamanik"
$ws.Range("D3").Value = 12000000
$ws.Range("E3").Value = "Manajemen Rekayasa"
$ws.Range("F3").Value = "PMDK"
$ws.Range("G3").Value = 3.3
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 3.4
$ws.Range("J3").Value = 3.5

# Row 4 new data
$ws.Range("A4").Value = 34567
$ws.Range("B4").Value = "Uzumymw"
$ws.Range("C4").Value = "SMAN 4 BINJAI"
$ws.Range("D4").Value = 2000000
$ws.Range("E4").Value = "Teknik Bioproses"
$ws.Range("F4").Value = "USM3"
$ws.Range("G4").Value = 1.3
$ws.Range("H4").Value = 1.2
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 2.1
$ws.Range("K4").Value = "Tidak tepat waktu"
